$d = $word.ActiveDocument
$CR = [char]13

# --------------------------------------------------------------------------
# Original layout (4 paragraphs):
#   1: "Ultimatum Game Instructions or Responders"          (title, unchanged)
#   2: bookmark paragraph (_GoBack), no visible text          (unchanged)
#   3: "In the following task, ... 0MU for that trial. "
#   4: "<tab>We are asking you to answer ... proposed division."
#
# Target layout: paragraphs 3 and 4 move to just after the title (each
# preceded/followed by blank paragraphs), the bookmark paragraph stays where
# it is, and a large new "Scenario 1/2/3" block is inserted right after the
# bookmark paragraph, followed by three trailing blank paragraphs.
# --------------------------------------------------------------------------

# Step 0: capture the original text of paragraphs 3 and 4 before any
# structural changes are made. Range.Text includes the trailing paragraph
# mark character, so strip it off.
$origP3Text = $d.Paragraphs.Item(3).Range.Text
$origP4Text = $d.Paragraphs.Item(4).Range.Text
$origP3Text = $origP3Text.Substring(0, $origP3Text.Length - 1)
$origP4Text = $origP4Text.Substring(0, $origP4Text.Length - 1)

# Step 1: insert the new "Scenario" content block right after the bookmark
# paragraph (paragraph 2). It consists of one leading blank paragraph, the
# scenario text (with blank-line separators), and three trailing blank
# paragraphs.
$scenarioBlock = (
  "Scenario 1" + $CR +
  "In this scenario you start the trial with 0 MU, and the proposer is given 20 MU to divide between the two of you. Therefore, the proposer has all of the money available to divide for this trial. In other words, the total amount of money at stake for the trial is 20 MU, and the proposer decides how to divide this." + $CR +
  "" + $CR +
  "For the explanations next to each offer:" + $CR +
  "The proposer offers you 0, so he or she receives 20, you receive 0. Etc." + $CR +
  "" + $CR +
  "Scenario 2" + $CR +
  "In this scenario you start each trial with 10 MU, and the proposer starts with 20 MU. Therefore, the total amount of money at stake for this trial is 30 MU, and the proposer decides how much of his 20 MU to divide between the two of you." + $CR +
  "" + $CR +
  "For the explanations next to each offer:" + $CR +
  "The proposer offers you 0, so he or she receives 20, and you receive 10." + $CR +
  "" + $CR +
  "Scenario 3" + $CR +
  "In this scenario you start each trial with 20 MU, and the proposer starts with 20 MU. Therefore, the total amount of money at stake for this trial is 40 MU, and the proposer decides how much of his 20 MU to divide between the two of you." + $CR +
  "" + $CR +
  "For the explanations next to each offer:" + $CR +
  "The proposer offers you 0, so he or she receives 20, and you receive 20." + $CR +
  "" + $CR +
  "" + $CR +
  ""
)

$r2 = $d.Paragraphs.Item(2).Range
$null = $r2.InsertParagraphAfter()   # becomes the blank paragraph before "Scenario 1"
$null = $r2.InsertParagraphAfter()   # becomes the paragraph that will hold the scenario text
$scenarioHost = $d.Paragraphs.Item(4)
$scenarioHost.Range.InsertBefore($scenarioBlock)

# Step 2: delete the original paragraphs 3 ("In the following...") and 4
# ("We are asking...") using Find, since at this point they are still the
# only occurrences of that text in the document (the duplicated copies for
# Block A have not been created yet).
$searchRange = $d.Content
$null = $searchRange.Find.Execute("In the following task,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$searchRange.Expand(4)
$searchRange.Delete()

$searchRange2 = $d.Content
$null = $searchRange2.Find.Execute("We are asking you to answer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$searchRange2.Expand(4)
$searchRange2.Delete()

# Step 3: insert Block A (blank, duplicated paragraph 3 text, duplicated
# paragraph 4 text, blank) right after the title paragraph (paragraph 1),
# i.e. just before the bookmark paragraph (still paragraph 2 at this point).
$blockA = (
  "" + $CR +
  $origP3Text + $CR +
  $origP4Text + $CR +
  ""
)
$r1 = $d.Paragraphs.Item(1).Range
$null = $r1.InsertParagraphAfter()
$blockAHost = $d.Paragraphs.Item(2)
$blockAHost.Range.InsertBefore($blockA)
